$d = $word.ActiveDocument

$pairs = @(
    ,@("15+48=", "47+25=")
    ,@("44+37=", "90-29=")
    ,@("9+4=", "33-4=")
    ,@("19+45=", "72-69=")
    ,@("50-32=", "80-32=")
    ,@("50-37=", "68+25=")
    ,@("70-69=", "17+28=")
    ,@("13+68=", "54+7=")
    ,@("80-64=", "66-27=")
    ,@("98-69=", "79+7=")
    ,@("9+24=", "36+37=")
    ,@("25+48=", "24+38=")
    ,@("28+7=", "37+26=")
    ,@("23+38=", "57-8=")
    ,@("70-63=", "58+23=")
    ,@("53+8=", "75+19=")
    ,@("72-56=", "29+19=")
    ,@("27+66=", "41-8=")
    ,@("73-67=", "83+8=")
    ,@("40-18=", "38+39=")
    ,@("71-17=", "4+57=")
    ,@("58+8=", "38+25=")
    ,@("6+29=", "71-66=")
    ,@("47+29=", "33+8=")
    ,@("80-9=", "19+26=")
    ,@("74-68=", "50-35=")
    ,@("90-18=", "83-78=")
    ,@("24+49=", "72-37=")
    ,@("91-89=", "19+73=")
    ,@("8+17=", "20-14=")
    ,@("51-47=", "80-24=")
    ,@("35+36=", "43-25=")
    ,@("77+19=", "57+26=")
    ,@("40-26=", "54+9=")
    ,@("8+14=", "29+55=")
    ,@("13+79=", "29+28=")
    ,@("30-18=", "7+34=")
    ,@("67-9=", "61-8=")
    ,@("8+13=", "55-29=")
    ,@("27+67=", "18+45=")
    ,@("17+26=", "90-82=")
    ,@("16+77=", "75-26=")
    ,@("97-39=", "71-3=")
    ,@("71-22=", "38+8=")
    ,@("79+8=", "87-49=")
    ,@("83-39=", "27+54=")
    ,@("37+54=", "56+28=")
    ,@("70-46=", "53-16=")
    ,@("8+24=", "55+39=")
    ,@("82-55=", "4+47=")
    ,@("26+49=", "7+69=")
    ,@("50-11=", "59+28=")
    ,@("39+14=", "44+28=")
    ,@("6+55=", "10-5=")
    ,@("92-5=", "46+49=")
    ,@("54-48=", "71-64=")
    ,@("47+49=", "82-43=")
    ,@("45+17=", "62-58=")
    ,@("92-17=", "22-17=")
    ,@("28+18=", "15-9=")
    ,@("48+6=", "86+5=")
    ,@("78+16=", "47+38=")
    ,@("66+5=", "84-8=")
    ,@("55-38=", "3+68=")
    ,@("38+34=", "10-5=")
    ,@("82-69=", "72-18=")
    ,@("27+36=", "80-77=")
    ,@("46-29=", "96-8=")
    ,@("19+68=", "18+47=")
    ,@("7+24=", "83-4=")
    ,@("61-12=", "8+75=")
    ,@("49+2=", "47+34=")
    ,@("25+39=", "85-49=")
    ,@("50-33=", "34-6=")
    ,@("84+8=", "74-19=")
    ,@("79+18=", "5+57=")
    ,@("88-29=", "78+9=")
    ,@("64-35=", "51-18=")
    ,@("61-17=", "41-3=")
    ,@("81-19=", "92-29=")
    ,@("46+25=", "25-7=")
    ,@("16+49=", "82-35=")
    ,@("68+23=", "60-25=")
    ,@("82-9=", "4+47=")
    ,@("27+29=", "60-32=")
    ,@("27-18=", "30-3=")
    ,@("87+8=", "37+36=")
    ,@("19+69=", "8+29=")
    ,@("43+9=", "54-7=")
    ,@("9+29=", "78+3=")
    ,@("30-19=", "29+48=")
    ,@("7+64=", "49+8=")
    ,@("71-29=", "51-34=")
    ,@("26-7=", "18+28=")
    ,@("98-49=", "38+43=")
    ,@("49+39=", "34-29=")
    ,@("88+7=", "14+19=")
    ,@("7+29=", "36+19=")
    ,@("29+62=", "88+8=")
    ,@("23-17=", "91-52=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
